$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new paragraph (two runs) right after the "How is code reuse
#    taken care of?" heading, before the following blank paragraph.
# ---------------------------------------------------------------------------
$findRng1 = $d.Content
[void]$findRng1.Find.Execute("How is code reuse taken care of?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target1 = $findRng1.Paragraphs(1)

$blank1 = $target1.Next()
$blank1.Range.InsertParagraphBefore()
$newPara1 = $target1.Next()

$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Code Reuse was something done reasonably well in this program, the major way was through inheritance, the majority of the class in the game inherit from the game object class, that contains variables and functions that get used in all of the children, the same also occurs with the projectile class, where projectiles of different types inherit from a parent projectile.</w:t></w:r><w:r><w:t xml:space="preserve"> It can also been seen in the AI where multiple AI can use the same state, in the case of our game the seek behaviour is used by two different AI' + [char]8217 + 's who simply set their own speed to differ from each other.</w:t></w:r></w:p>'
[void]$newPara1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Replace the final paragraph: drop the old "How would this reuse
#    happen?" run, keep the _GoBack bookmark, and append two new runs.
# ---------------------------------------------------------------------------
$findRng2 = $d.Content
[void]$findRng2.Find.Execute("How would this reuse happen?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target2 = $findRng2.Paragraphs(1)

$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Some of the parent classes could be very useful in other programs, the gameobject class is a basic class that contains variables and functions of use to any object placed in a game, this would be quite useful in another game program however details to do with displaying may need to be changed.</w:t></w:r><w:r><w:t xml:space="preserve"> The finite state machine and associated state class are probably the most portable pieces of code, zero changes would be required to use these classes in another program and it would not be restricted to game. The same also holds true for the sound class.</w:t></w:r></w:p>'
[void]$target2.Range.InsertXML($xml2)
